$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records to append below the existing data (update through 02/05).
$newRows = @(
    @(44313, 2, 17, 133.4903808402042),
    @(44314, 4, 21, 164.8998822143699),
    @(44315, 4, 22, 172.7522575579113),
    @(44316, 5, 26, 204.161758932077),
    @(44317, 5, 30, 235.5712603062426),
    @(44318, 0, 25, 196.3093835885355)
)

$lastRow = 238
$startRow = $lastRow + 1

$r = $startRow
foreach ($rowData in $newRows) {
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]

    # Match the date-column formatting (s="2") used by the rest of column A.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $r++
}

$excel.CutCopyMode = 0

Write-Output "Appended rows $startRow to $($r - 1)"
